$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-10 down to 8-11
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the new weekly entry
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44839
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 26000
$ws.Range("P7").Value = 25500
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 2550
$ws.Range("T7").Value = 10
